$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.088.79"
$ws.Range("E2").Value = "  +1.53%  "
$ws.Range("D3").Value = "2.508.93"
$ws.Range("E3").Value = "  +0.90%  "
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("D5").Value = "'536.80"
$ws.Range("E5").Value = "  +3.42%  "
$ws.Range("D6").Value = "'134.83"
$ws.Range("E6").Value = "  +2.26%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("E8").Value = "  +2.98%  "
$ws.Range("D9").Value = "2.514.45"
$ws.Range("E9").Value = "  -0.29%  "
$ws.Range("D10").Value = "'0.100"
$ws.Range("E10").Value = "  +3.19%  "
$ws.Range("E11").Value = "  -2.76%  "
$ws.Range("D12").Value = "'5.18"
$ws.Range("E12").Value = "  -0.25%  "
$ws.Range("E13").Value = "  -0.47%  "
$ws.Range("D14").Value = "2.951.64"
$ws.Range("E14").Value = "  -0.42%  "
$ws.Range("D15").Value = "58.860.30"
$ws.Range("E15").Value = "  +1.19%  "
$ws.Range("D16").Value = "'22.44"
$ws.Range("E16").Value = "  +1.48%  "
$ws.Range("E17").Value = "  +1.12%  "
$ws.Range("D18").Value = "2.513.08"
$ws.Range("E18").Value = "  -0.02%  "
$ws.Range("D19").Value = "'10.71"
$ws.Range("E19").Value = "  +0.07%  "
$ws.Range("D20").Value = "'4.26"
$ws.Range("E20").Value = "  +2.01%  "
$ws.Range("D21").Value = "'321.91"
$ws.Range("E21").Value = "  +0.20%  "
$ws.Range("E22").Value = "  +4.66%  "
$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("E24").Value = "  +2.35%  "
$ws.Range("E25").Value = "  +1.76%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.45%  "
$ws.Range("E27").Value = "  -0.15%  "
$ws.Range("E28").Value = "  +1.83%  "
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").Value = "'174.13"
$ws.Range("E29").Value = "  +3.64%  "
$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").Value = "0.0₃0762"
$ws.Range("E30").Value = "  +1.83%  "
$ws.Range("D31").Value = "'1.74"
$ws.Range("D32").Value = "'1.19"
$ws.Range("E32").Value = "  +1.59%  "
$ws.Range("D33").Value = "'6.29"
$ws.Range("E33").Value = "  +0.85%  "
$ws.Range("E34").Value = "  +0.07%  "
$ws.Range("E35").Value = "  +0.36%  "
$ws.Range("D36").Value = "'18.16"
$ws.Range("E36").Value = "  +0.54%  "
$ws.Range("E37").Value = "  -2.04%  "
$ws.Range("D38").Value = "'3.95"
$ws.Range("E38").Value = "  +0.40%  "
$ws.Range("E39").Value = "  +4.30%  "
$ws.Range("D40").Value = "'0.826"
$ws.Range("E40").Value = "  +7.48%  "
$ws.Range("D41").Value = "'36.63"
$ws.Range("E41").Value = "  -0.89%  "
$ws.Range("D42").Value = "'3.50"
$ws.Range("E42").Value = "  +1.73%  "
$ws.Range("D43").Value = "'276.78"
$ws.Range("E43").Value = "  -0.04%  "
$ws.Range("D44").Value = "'132.01"
$ws.Range("E44").Value = "  +8.60%  "
$ws.Range("E45").Value = "  -0.70%  "
$ws.Range("D46").Value = "'0.593"
$ws.Range("E46").Value = "  -0.70%  "
$ws.Range("D47").Value = "'0.0944"
$ws.Range("E47").Value = "  +2.60%  "
$ws.Range("E48").Value = "  +2.51%  "
$ws.Range("D49").Value = "'0.0220"
$ws.Range("E49").Value = "  +3.22%  "
$ws.Range("D50").Value = "'16.94"
$ws.Range("E50").Value = "  +0.21%  "
$ws.Range("D51").Value = "1.755.02"
$ws.Range("E51").Value = "  +0.93%  "
